$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new label "nível:" in D4 (same style/format as other label cells in column A, i.e. s=4)
$ws.Range("A4").Copy()
$ws.Range("D4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D4").Value = "nível:"

# Add new header "Observações" in O6 (same style as header row, s=1 -> bordered, centered, bold Arial 10)
# Copy the exact formatting already used by the neighboring header cells (E6:N6)
$ws.Range("E6").Copy()
$ws.Range("O6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("O6").Value = "Observações"

# Column O width to fit new content (approximates Excel's best-fit autosize result)
$ws.Range("O1").EntireColumn.ColumnWidth = 11.736979166666666

# Update selection to D2 to mirror the saved view state
$ws.Range("D2").Select()
